$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 18
$ws_ALC.Range("H18").Value = 83334840
$ws_ALC.Range("I18").Value = 111112450
$ws_ALC.Range("K18").Value = 111112450
$ws_ALC.Range("M18").Value = -111112166

# ALC row 19
$ws_ALC.Range("H19").Value = 991.6667
$ws_ALC.Range("I19").Value = 999.5
$ws_ALC.Range("J19").Value = 989.4286
$ws_ALC.Range("K19").Value = 999.5
$ws_ALC.Range("L19").Value = 989.4286
$ws_ALC.Range("M19").Value = -824.5
$ws_ALC.Range("N19").Value = -1339.4286

# ALC row 40
$ws_ALC.Range("H40").Value = 11380658
$ws_ALC.Range("I40").Value = 9554.964
$ws_ALC.Range("J40").Value = 31280088
$ws_ALC.Range("K40").Value = 9554.964
$ws_ALC.Range("L40").Value = 31280088
$ws_ALC.Range("M40").Value = -9379.964
$ws_ALC.Range("N40").Value = -31280438

# ALC row 80
$ws_ALC.Range("H80").Value = 1058.7
$ws_ALC.Range("I80").Value = 755.2857
$ws_ALC.Range("K80").Value = 2265.8571
$ws_ALC.Range("M80").Value = -1267.8571

# ALC row 83
$ws_ALC.Range("H83").Value = 1058.7
$ws_ALC.Range("I83").Value = 755.2857
$ws_ALC.Range("K83").Value = 6797.571300000001
$ws_ALC.Range("M83").Value = -1805.571300000001

# ALC row 96
$ws_ALC.Range("H96").Value = 18788.445
$ws_ALC.Range("I96").Value = 737
$ws_ALC.Range("J96").Value = 33229.6
$ws_ALC.Range("K96").Value = 2211
$ws_ALC.Range("L96").Value = 99688.79999999999
$ws_ALC.Range("M96").Value = -838
$ws_ALC.Range("N96").Value = -102434.8

# ALC row 100
$ws_ALC.Range("H100").Value = 2441.5
$ws_ALC.Range("I100").Value = 2441.5
$ws_ALC.Range("K100").Value = 2441.5
$ws_ALC.Range("M100").Value = -1900.5

# ALC row 140
$ws_ALC.Range("H140").Value = 62836.445
$ws_ALC.Range("J140").Value = 61852.375
$ws_ALC.Range("L140").Value = 61852.375
$ws_ALC.Range("N140").Value = -72212.375

# ARM row 32
$ws_ARM.Range("H32").Value = 3649.0986
$ws_ARM.Range("I32").Value = 2250.4666
$ws_ARM.Range("K32").Value = 2250.4666
$ws_ARM.Range("M32").Value = -1963.4666

# BSM row 134
$ws_BSM.Range("H134").Value = 1941.4615
$ws_BSM.Range("I134").Value = 691.4666999999999
$ws_BSM.Range("K134").Value = 2074.4001
$ws_BSM.Range("M134").Value = 460.5999000000002

# CRP row 105
$ws_CRP.Range("H105").Value = 1894970.9
$ws_CRP.Range("I105").Value = 2526050
$ws_CRP.Range("K105").Value = 2526050
$ws_CRP.Range("M105").Value = -2524303

# CRP row 124
$ws_CRP.Range("H124").Value = 44325.5
$ws_CRP.Range("J124").Value = 44325.5
$ws_CRP.Range("L124").Value = 44325.5
$ws_CRP.Range("N124").Value = -49235.5

# CUL row 4
$ws_CUL.Range("H4").Value = 5995790
$ws_CUL.Range("I4").Value = 2241870.2
$ws_CUL.Range("K4").Value = 6725610.600000001
$ws_CUL.Range("M4").Value = -6725498.600000001

# CUL row 19
$ws_CUL.Range("H19").Value = 100
$ws_CUL.Range("I19").Value = 100
$ws_CUL.Range("K19").Value = 300
$ws_CUL.Range("M19").Value = -126

# CUL row 68
$ws_CUL.Range("H68").Value = 1557105.4
$ws_CUL.Range("I68").Value = 1750
$ws_CUL.Range("J68").Value = 1945944.2
$ws_CUL.Range("K68").Value = 5250
$ws_CUL.Range("L68").Value = 5837832.6
$ws_CUL.Range("M68").Value = -4439
$ws_CUL.Range("N68").Value = -5839454.6

# CUL row 69
$ws_CUL.Range("H69").Value = 3661.4614
$ws_CUL.Range("I69").Value = 0
$ws_CUL.Range("J69").Value = 3661.4614
$ws_CUL.Range("K69").Value = 0
$ws_CUL.Range("L69").Value = 10984.3842
$ws_CUL.Range("M69").ClearContents()
$ws_CUL.Range("N69").Value = -12606.3842

# CUL row 71
$ws_CUL.Range("H71").Value = 1557105.4
$ws_CUL.Range("I71").Value = 1750
$ws_CUL.Range("J71").Value = 1945944.2
$ws_CUL.Range("K71").Value = 15750
$ws_CUL.Range("L71").Value = 17513497.8
$ws_CUL.Range("M71").Value = -11694
$ws_CUL.Range("N71").Value = -17521609.8

# CUL row 72
$ws_CUL.Range("H72").Value = 3661.4614
$ws_CUL.Range("I72").Value = 0
$ws_CUL.Range("J72").Value = 3661.4614
$ws_CUL.Range("K72").Value = 0
$ws_CUL.Range("L72").Value = 32953.1526
$ws_CUL.Range("M72").ClearContents()
$ws_CUL.Range("N72").Value = -41065.1526

# CUL row 80
$ws_CUL.Range("H80").Value = 5999.125
$ws_CUL.Range("I80").Value = 0
$ws_CUL.Range("J80").Value = 5999.125
$ws_CUL.Range("K80").Value = 0
$ws_CUL.Range("L80").Value = 17997.375
$ws_CUL.Range("M80").ClearContents()
$ws_CUL.Range("N80").Value = -19869.375

# CUL row 83
$ws_CUL.Range("H83").Value = 5999.125
$ws_CUL.Range("I83").Value = 0
$ws_CUL.Range("J83").Value = 5999.125
$ws_CUL.Range("K83").Value = 0
$ws_CUL.Range("L83").Value = 53992.125
$ws_CUL.Range("M83").ClearContents()
$ws_CUL.Range("N83").Value = -63352.125

# CUL row 86
$ws_CUL.Range("H86").Value = 281.25
$ws_CUL.Range("I86").Value = 281.25
$ws_CUL.Range("J86").Value = 0
$ws_CUL.Range("K86").Value = 843.75
$ws_CUL.Range("L86").Value = 0
$ws_CUL.Range("M86").Value = 342.25
$ws_CUL.Range("N86").ClearContents()

# CUL row 89
$ws_CUL.Range("H89").Value = 281.25
$ws_CUL.Range("I89").Value = 281.25
$ws_CUL.Range("J89").Value = 0
$ws_CUL.Range("K89").Value = 2531.25
$ws_CUL.Range("L89").Value = 0
$ws_CUL.Range("M89").Value = 3396.75
$ws_CUL.Range("N89").ClearContents()

# GSM row 43
$ws_GSM.Range("H43").Value = 21799.8
$ws_GSM.Range("I43").Value = 17999.666
$ws_GSM.Range("K43").Value = 17999.666
$ws_GSM.Range("M43").Value = -17848.666

# GSM row 52
$ws_GSM.Range("H52").Value = 12000
$ws_GSM.Range("I52").Value = 12000
$ws_GSM.Range("K52").Value = 12000
$ws_GSM.Range("M52").Value = -11741

# GSM row 113
$ws_GSM.Range("H113").Value = 1557.8334
$ws_GSM.Range("I113").Value = 1424.5
$ws_GSM.Range("J113").Value = 1584.5
$ws_GSM.Range("K113").Value = 1424.5
$ws_GSM.Range("L113").Value = 1584.5
$ws_GSM.Range("M113").Value = 745.5
$ws_GSM.Range("N113").Value = -5924.5

# LTW row 61
$ws_LTW.Range("H61").Value = 1910
$ws_LTW.Range("I61").Value = 1581.762
$ws_LTW.Range("K61").Value = 1581.762
$ws_LTW.Range("M61").Value = -1379.762

# LTW row 100
$ws_LTW.Range("H100").Value = 5092.6665
$ws_LTW.Range("I100").Value = 5450
$ws_LTW.Range("J100").Value = 4914
$ws_LTW.Range("K100").Value = 5450
$ws_LTW.Range("L100").Value = 4914
$ws_LTW.Range("M100").Value = -4909
$ws_LTW.Range("N100").Value = -5996

# LTW row 113
$ws_LTW.Range("H113").Value = 1910
$ws_LTW.Range("I113").Value = 1581.762
$ws_LTW.Range("K113").Value = 1581.762
$ws_LTW.Range("M113").Value = 588.2380000000001

# LTW row 122
$ws_LTW.Range("H122").Value = 57146710
$ws_LTW.Range("I122").Value = 125003540
$ws_LTW.Range("K122").Value = 375010620
$ws_LTW.Range("M122").Value = -375008170

# WVR row 38
$ws_WVR.Range("H38").Value = 12000
$ws_WVR.Range("J38").Value = 12000
$ws_WVR.Range("L38").Value = 12000
$ws_WVR.Range("N38").Value = -12946

# WVR row 49
$ws_WVR.Range("H49").Value = 65000
$ws_WVR.Range("J49").Value = 65000
$ws_WVR.Range("L49").Value = 65000
$ws_WVR.Range("N49").Value = -65460

# WVR row 96
$ws_WVR.Range("H96").Value = 5362.6
$ws_WVR.Range("J96").Value = 5453.5
$ws_WVR.Range("L96").Value = 5453.5
$ws_WVR.Range("N96").Value = -8199.5
